$wb = $excel.ActiveWorkbook

# Overview sheet: row 3 (c160a3b9 file) status moves from "Ready for handoff"
# to "Handed back: in sync with en-US"
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C3").Value = "Handed back: in sync with en-US"

# zh-cn sheet: same status update, plus new Latest Handback DateTime
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("B3").Value = "Handed back: in sync with en-US"
$wsZh.Range("G3").Value = "2016-01-11 16:39:48"

# de-de sheet: same status update, plus new Latest Handback DateTime
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("B3").Value = "Handed back: in sync with en-US"
$wsDe.Range("G3").Value = "2016-01-11 16:40:27"
